$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Rendi Imam Saputra"),
    @(1, "Dede Ahmad Fauzan"),
    @(1, "Wilda Nazwatun Nisa"),
    @(2, "Hoirul Sambudi"),
    @(2, "Dinda Ayuni"),
    @(2, "Ardan Mizanul Khoiri"),
    @(3, "Agung Prayuda"),
    @(3, "Citra Hayatunnufus"),
    @(3, "Mochammad Wafi Nur Jihan"),
    @(4, "Khairun Anwar"),
    @(4, "Muhamad Maulana Naufal Widodo"),
    @(4, "Ramlan"),
    @(5, "Nurul Hidayah Harahap"),
    @(5, "Dwi Nur Aini"),
    @(5, "Fadli Al Masani")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 2).Style = "Normal"
}

$ws.Range("D8").Select()
